$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.8639936447143555
$ws.Range("E2").Value = 296.1195148364441
$ws.Range("F2").Value = 0.01065921349924351
$ws.Range("G2").Value = 0.009327086568022135
$ws.Range("H2").Value = 0.008519410935977854
$ws.Range("I2").Value = 0.008051945088925677
$ws.Range("J2").Value = 0.007303639323851956
$ws.Range("K2").Value = 0.007141400872617726
$ws.Range("L2").Value = 0.006779025560253142
$ws.Range("M2").Value = 0.006779025560253142
$ws.Range("N2").Value = 0.006561313563694094
$ws.Range("O2").Value = 0.006397935445406792
$ws.Range("P2").Value = 0.006266892958815288
$ws.Range("Q2").Value = 0.00618954644162144
$ws.Range("R2").Value = 0.00603937521859977
$ws.Range("S2").Value = 0.005976829282653798
$ws.Range("T2").Value = 0.005947638266956749
$ws.Range("U2").Value = 0.005907942566718205
$ws.Range("V2").Value = 0.005859474980700902
$ws.Range("W2").Value = 0.005824599767077381
$ws.Range("X2").Value = 0.005791194632634606
$ws.Range("Y2").Value = 0.005772310230729902

$ws.Range("C3").Value = 0.986997127532959
$ws.Range("E3").Value = 289.2158871890515
$ws.Range("F3").Value = 0.01093838845198985
$ws.Range("G3").Value = 0.009110775181622359
$ws.Range("H3").Value = 0.008089553849284475
$ws.Range("I3").Value = 0.007563562145410916
$ws.Range("J3").Value = 0.007321821088730339
$ws.Range("K3").Value = 0.007172755983415963
$ws.Range("L3").Value = 0.006745163143400894
$ws.Range("M3").Value = 0.006561719699221056
$ws.Range("N3").Value = 0.00639307252496842
$ws.Range("O3").Value = 0.006154336109788412
$ws.Range("P3").Value = 0.006102585724693127
$ws.Range("Q3").Value = 0.005972710708785497
$ws.Range("R3").Value = 0.005880166645304227
$ws.Range("S3").Value = 0.005839371578010052
$ws.Range("T3").Value = 0.005770312520486029
$ws.Range("U3").Value = 0.005751366960597891
$ws.Range("V3").Value = 0.0057095001965325
$ws.Range("W3").Value = 0.005678758001815071
$ws.Range("X3").Value = 0.005657860847128233
$ws.Range("Y3").Value = 0.005637736592379171

$ws.Range("C4").Value = 0.8109967708587646
$ws.Range("E4").Value = 281.3937139852096
$ws.Range("F4").Value = 0.0109174859171614
$ws.Range("G4").Value = 0.009000449033074084
$ws.Range("H4").Value = 0.0079721850354861
$ws.Range("I4").Value = 0.007524022928289915
$ws.Range("J4").Value = 0.00713000066991834
$ws.Range("K4").Value = 0.006729976223824906
$ws.Range("L4").Value = 0.006330224090856857
$ws.Range("M4").Value = 0.006114075027566264
$ws.Range("N4").Value = 0.005967455159733072
$ws.Range("O4").Value = 0.005896934349733619
$ws.Range("P4").Value = 0.005839350711776134
$ws.Range("Q4").Value = 0.005773278928779519
$ws.Range("R4").Value = 0.005746901577213047
$ws.Range("S4").Value = 0.005672201782516224
$ws.Range("T4").Value = 0.005606174398824691
$ws.Range("U4").Value = 0.005570916852754136
$ws.Range("V4").Value = 0.005524872040737942
$ws.Range("W4").Value = 0.00550516647454677
$ws.Range("X4").Value = 0.005488800672767955
$ws.Range("Y4").Value = 0.005485257582557691

$ws.Range("C5").Value = 0.8500010967254639
$ws.Range("E5").Value = 287.8701452900732
$ws.Range("F5").Value = 0.01082266316085568
$ws.Range("G5").Value = 0.009187986263855315
$ws.Range("H5").Value = 0.008278100275909318
$ws.Range("I5").Value = 0.007433804296935851
$ws.Range("J5").Value = 0.006936204093651316
$ws.Range("K5").Value = 0.006853270744947282
$ws.Range("L5").Value = 0.006484330964937192
$ws.Range("M5").Value = 0.006397406068918979
$ws.Range("N5").Value = 0.006249074046670152
$ws.Range("O5").Value = 0.006166039595275993
$ws.Range("P5").Value = 0.006012928734474868
$ws.Range("Q5").Value = 0.005948819323934726
$ws.Range("R5").Value = 0.005857873041175507
$ws.Range("S5").Value = 0.005808114280872779
$ws.Range("T5").Value = 0.005739185838940374
$ws.Range("U5").Value = 0.005715930484355331
$ws.Range("V5").Value = 0.005698095070917781
$ws.Range("W5").Value = 0.005642685121888009
$ws.Range("X5").Value = 0.005630724235496577
$ws.Range("Y5").Value = 0.005611503806824038

$ws.Range("C6").Value = 0.8560261726379395
$ws.Range("E6").Value = 287.9387669249809
$ws.Range("F6").Value = 0.01079618492857642
$ws.Range("G6").Value = 0.009222670275773202
$ws.Range("H6").Value = 0.008485776106627934
$ws.Range("I6").Value = 0.007701904731073876
$ws.Range("J6").Value = 0.007284695299467978
$ws.Range("K6").Value = 0.007040250758469538
$ws.Range("L6").Value = 0.006999540033925368
$ws.Range("M6").Value = 0.006649806100610757
$ws.Range("N6").Value = 0.00639928383211524
$ws.Range("O6").Value = 0.00622351316853138
$ws.Range("P6").Value = 0.006059245119800802
$ws.Range("Q6").Value = 0.005990242035827284
$ws.Range("R6").Value = 0.005892514446745176
$ws.Range("S6").Value = 0.005775424016101026
$ws.Range("T6").Value = 0.005775424016101026
$ws.Range("U6").Value = 0.005669662954929365
$ws.Range("V6").Value = 0.005669379922316247
$ws.Range("W6").Value = 0.005659077621671153
$ws.Range("X6").Value = 0.005612841460525942
$ws.Range("Y6").Value = 0.005612841460525942

$ws.Range("C7").Value = 0.7842490673065186
$ws.Range("E7").Value = 281.022436977104
$ws.Range("F7").Value = 0.01073066322429305
$ws.Range("G7").Value = 0.009021469864119475
$ws.Range("H7").Value = 0.008242146319068928
$ws.Range("I7").Value = 0.007767383393101018
$ws.Range("J7").Value = 0.007115233313051412
$ws.Range("K7").Value = 0.006825990454005032
$ws.Range("L7").Value = 0.006543071547311605
$ws.Range("M7").Value = 0.006231704322282084
$ws.Range("N7").Value = 0.006107595443244023
$ws.Range("O7").Value = 0.006047382643343071
$ws.Range("P7").Value = 0.005961716176295798
$ws.Range("Q7").Value = 0.005812304131772955
$ws.Range("R7").Value = 0.005812304131772955
$ws.Range("S7").Value = 0.005732510864528489
$ws.Range("T7").Value = 0.00564571010558406
$ws.Range("U7").Value = 0.005599221026200645
$ws.Range("V7").Value = 0.005553288670629906
$ws.Range("W7").Value = 0.005532114010930156
$ws.Range("X7").Value = 0.005498924641649697
$ws.Range("Y7").Value = 0.005478020213978634

$ws.Range("C8").Value = 0.6970038414001465
$ws.Range("E8").Value = 287.2063047151532
$ws.Range("F8").Value = 0.01069053197663875
$ws.Range("G8").Value = 0.009258540707633107
$ws.Range("H8").Value = 0.00841835531749475
$ws.Range("I8").Value = 0.007803719938648874
$ws.Range("J8").Value = 0.007190188045008193
$ws.Range("K8").Value = 0.007018750288422422
$ws.Range("L8").Value = 0.006825461038254706
$ws.Range("M8").Value = 0.006635392794623432
$ws.Range("N8").Value = 0.006215421924813485
$ws.Range("O8").Value = 0.006215421924813485
$ws.Range("P8").Value = 0.006112963457219991
$ws.Range("Q8").Value = 0.005988202346238889
$ws.Range("R8").Value = 0.005903451555463262
$ws.Range("S8").Value = 0.005812518720200638
$ws.Range("T8").Value = 0.005675801240050402
$ws.Range("U8").Value = 0.005675801240050402
$ws.Range("V8").Value = 0.005675801240050402
$ws.Range("W8").Value = 0.005635402543683981
$ws.Range("X8").Value = 0.005627217927530435
$ws.Range("Y8").Value = 0.005598563444739828

$ws.Range("C9").Value = 0.7299864292144775
$ws.Range("E9").Value = 300.1896245735297
$ws.Range("F9").Value = 0.01074539905723561
$ws.Range("G9").Value = 0.009039369785125803
$ws.Range("H9").Value = 0.008300985575273402
$ws.Range("I9").Value = 0.007710844899008651
$ws.Range("J9").Value = 0.007247188991277652
$ws.Range("K9").Value = 0.006855724593098206
$ws.Range("L9").Value = 0.006676931233639734
$ws.Range("M9").Value = 0.006476097353256385
$ws.Range("N9").Value = 0.006329770271015865
$ws.Range("O9").Value = 0.006255902372528957
$ws.Range("P9").Value = 0.006201293997590424
$ws.Range("Q9").Value = 0.006147373706760837
$ws.Range("R9").Value = 0.006042388022378269
$ws.Range("S9").Value = 0.005999455501715506
$ws.Range("T9").Value = 0.005969556808319826
$ws.Range("U9").Value = 0.00594910304184398
$ws.Range("V9").Value = 0.005915342941916073
$ws.Range("W9").Value = 0.005888226746343741
$ws.Range("X9").Value = 0.005871895406501307
$ws.Range("Y9").Value = 0.005851649601823191

$ws.Range("C10").Value = 0.9390068054199219
$ws.Range("E10").Value = 292.0559035863625
$ws.Range("F10").Value = 0.01052853774516058
$ws.Range("G10").Value = 0.009174170963945861
$ws.Range("H10").Value = 0.008376687239522661
$ws.Range("I10").Value = 0.007462415807519122
$ws.Range("J10").Value = 0.007127028139438352
$ws.Range("K10").Value = 0.006887780987552167
$ws.Range("L10").Value = 0.00680813666965815
$ws.Range("M10").Value = 0.006410018289955128
$ws.Range("N10").Value = 0.00635326265732845
$ws.Range("O10").Value = 0.006243296431638505
$ws.Range("P10").Value = 0.00614398900002081
$ws.Range("Q10").Value = 0.006010514260017062
$ws.Range("R10").Value = 0.00589372902751667
$ws.Range("S10").Value = 0.00589372902751667
$ws.Range("T10").Value = 0.005835714732511688
$ws.Range("U10").Value = 0.005773394274065121
$ws.Range("V10").Value = 0.005740561564248813
$ws.Range("W10").Value = 0.005706293366762556
$ws.Range("X10").Value = 0.005706293366762556
$ws.Range("Y10").Value = 0.005693097535796539

$ws.Range("C11").Value = 0.7250375747680664
$ws.Range("E11").Value = 289.8192262131215
$ws.Range("F11").Value = 0.01081722240702897
$ws.Range("G11").Value = 0.009084748934296262
$ws.Range("H11").Value = 0.008245047107720962
$ws.Range("I11").Value = 0.007796060707096364
$ws.Range("J11").Value = 0.007490931808848297
$ws.Range("K11").Value = 0.007268194483581432
$ws.Range("L11").Value = 0.006907612672054565
$ws.Range("M11").Value = 0.006612452277113794
$ws.Range("N11").Value = 0.006496602793851128
$ws.Range("O11").Value = 0.00605310384628778
$ws.Range("P11").Value = 0.00605310384628778
$ws.Range("Q11").Value = 0.005977578283075318
$ws.Range("R11").Value = 0.005949941000725935
$ws.Range("S11").Value = 0.005893052227192949
$ws.Range("T11").Value = 0.005752749214297771
$ws.Range("U11").Value = 0.005752749214297771
$ws.Range("V11").Value = 0.005752749214297771
$ws.Range("W11").Value = 0.005709073737954795
$ws.Range("X11").Value = 0.005674341847459206
$ws.Range("Y11").Value = 0.005649497587000418
